$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts, previously Strike#) values recomputed for rows 2..34.
$newK = @{
    2  = 2
    3  = 0
    4  = 2
    5  = 1
    6  = 3
    7  = 1
    8  = 2
    9  = 4
    10 = 0
    11 = 3
    12 = 2
    13 = 3
    14 = 1
    15 = 0
    16 = 1
    17 = 3
    18 = 3
    19 = 1
    20 = 4
    21 = 3
    22 = 4
    23 = 2
    24 = 2
    25 = 8
    26 = 1
    27 = 6
    28 = 2
    29 = 2
    30 = 1
    31 = 3
    32 = 2
    33 = 1
    34 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
